# Recalculated astrometric residual columns (ra_resid_deg, dec_resid_deg,
# ra_resid_rad, dec_resid_rad) for rows 2-21 using the updated variance
# calculation. Values below are the recomputed results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = 0.0001180168397070247
$ws.Range("AG2").Value = 0.000007043030873177258
$ws.Range("AH2").Value = 0.000002059782425685961
$ws.Range("AI2").Value = 0.0000001229240780565543

$ws.Range("AF3").Value = -0.00009763669208950887
$ws.Range("AG3").Value = 0.00002042721613548792
$ws.Range("AH3").Value = -0.000001704081747717832
$ws.Range("AI3").Value = 0.0000003565221785807764

$ws.Range("AF4").Value = 0.0005529119226252988
$ws.Range("AG4").Value = 0.0001569640776146741
$ws.Range("AH4").Value = 0.000009650133523343592
$ws.Range("AI4").Value = 0.00000273953996173199

$ws.Range("AF5").Value = -0.0002319378970128128
$ws.Range("AG5").Value = -0.00001761794495180879
$ws.Range("AH5").Value = -0.000004048079963025103
$ws.Range("AI5").Value = -0.0000003074911468441772

$ws.Range("AF6").Value = 0.0001555432800444123
$ws.Range("AG6").Value = 0.00001860249967577943
$ws.Range("AH6").Value = 0.000002714742366126587
$ws.Range("AI6").Value = 0.0000003246748684435286

$ws.Range("AF7").Value = 0.000001406725147035104
$ws.Range("AG7").Value = -0.000006499588135255863
$ws.Range("AH7").Value = 0.00000002455198548636391
$ws.Range("AI7").Value = -0.0000001134392129837733

$ws.Range("AF8").Value = 0.00004557408715299971
$ws.Range("AG8").Value = 0.00001121169369966424
$ws.Range("AH8").Value = 0.0000007954178744106936
$ws.Range("AI8").Value = 0.0000001956809697842452

$ws.Range("AF9").Value = 0.00006506150984364467
$ws.Range("AG9").Value = 0.0000007226120040115802
$ws.Range("AH9").Value = 0.000001135537563090301
$ws.Range("AI9").Value = 0.00000001261195868443655

$ws.Range("AF10").Value = -0.00008563498568037176
$ws.Range("AG10").Value = -0.000006624645177311095
$ws.Range("AH10").Value = -0.000001494612455020684
$ws.Range("AI10").Value = -0.0000001156218701204422

$ws.Range("AF11").Value = 0.00009742706333781825
$ws.Range("AG11").Value = 0.00002122749405941704
$ws.Range("AH11").Value = 0.000001700423035793985
$ws.Range("AI11").Value = 0.0000003704896632843641

$ws.Range("AF12").Value = -0.0001139440585689044
$ws.Range("AG12").Value = -0.00006762308426644381
$ws.Range("AH12").Value = -0.00000198869898511264
$ws.Range("AI12").Value = -0.000001180245470803019

$ws.Range("AF13").Value = -0.00006190919964410568
$ws.Range("AG13").Value = -0.000001899857210219125
$ws.Range("AH13").Value = -0.000001080519371064146
$ws.Range("AI13").Value = -0.00000003315876363607778

$ws.Range("AF14").Value = -0.00005394209580344977
$ws.Range("AG14").Value = 0.00002020891774634492
$ws.Range("AH14").Value = -0.0000009414671771964145
$ws.Range("AI14").Value = 0.0000003527121529384311

$ws.Range("AF15").Value = -0.00009637072705004357
$ws.Range("AG15").Value = 0.00009357044092439537
$ws.Range("AH15").Value = -0.000001681986489564022
$ws.Range("AI15").Value = 0.000001633112276673546

$ws.Range("AF16").Value = 0.00008870599910437704
$ws.Range("AG16").Value = -0.000007537564876258784
$ws.Range("AH16").Value = 0.000001548211750642521
$ws.Range("AI16").Value = -0.0000001315553246733947

$ws.Range("AF17").Value = -0.0002119636181134865
$ws.Range("AG17").Value = -0.0002107375396747102
$ws.Range("AH17").Value = -0.000003699463030520231
$ws.Range("AI17").Value = -0.000003678063924875873

$ws.Range("AF18").Value = -0.00008169830471160822
$ws.Range("AG18").Value = -0.00003105098020128594
$ws.Range("AH18").Value = -0.00000142590441051516
$ws.Range("AI18").Value = -0.0000005419418404840111

$ws.Range("AF19").Value = -0.0001192559604135113
$ws.Range("AG19").Value = 0.00001236068153787073
$ws.Range("AH19").Value = -0.00000208140916184379
$ws.Range("AI19").Value = 0.0000002157345906263204

$ws.Range("AF20").Value = -0.00005736948804724307
$ws.Range("AG20").Value = 0.0000180555296438456
$ws.Range("AH20").Value = -0.000001001286456607924
$ws.Range("AI20").Value = 0.0000003151284404765449

$ws.Range("AF21").Value = 0.00008701560065560443
$ws.Range("AG21").Value = -0.00003080298937518933
$ws.Range("AH21").Value = 0.000001518708732040834
$ws.Range("AI21").Value = -0.0000005376135840538846
